$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 177.5
$ws.Range("I6").Value = 196.66667
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 590.00001
$ws.Range("L6").Value = 15
$ws.Range("M6").Value = -478.00001
$ws.Range("N6").Value = -239
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830
$ws.Range("H31").Value = 2666.3333
$ws.Range("I31").Value = 2666.3333
$ws.Range("K31").Value = 7998.999899999999
$ws.Range("M31").Value = -7768.999899999999
$ws.Range("H41").Value = 466.33334
$ws.Range("I41").Value = 494
$ws.Range("J41").Value = 328
$ws.Range("K41").Value = 494
$ws.Range("L41").Value = 328
$ws.Range("M41").Value = -54
$ws.Range("N41").Value = -1208
$ws.Range("H43").Value = 974.8
$ws.Range("I43").Value = 992.3333
$ws.Range("K43").Value = 992.3333
$ws.Range("M43").Value = -923.3333
$ws.Range("H70").Value = 2083
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 2124.5
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 6373.5
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -6913.5
$ws.Range("H73").Value = 2083
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 2124.5
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 6373.5
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -8245.5
$ws.Range("H96").Value = 10992.3
$ws.Range("I96").Value = 13053.125
$ws.Range("K96").Value = 39159.375
$ws.Range("M96").Value = -37786.375
$ws.Range("H125").Value = 2999.5
$ws.Range("J125").Value = 4999
$ws.Range("L125").Value = 44991
$ws.Range("N125").Value = -49911

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2854762.2
$ws.Range("J32").Value = 3500497
$ws.Range("L32").Value = 3500497
$ws.Range("N32").Value = -3501071
$ws.Range("H45").Value = 7077.5
$ws.Range("I45").Value = 9066.666999999999
$ws.Range("K45").Value = 9066.666999999999
$ws.Range("M45").Value = -8689.666999999999
$ws.Range("H122").Value = 2326.7
$ws.Range("I122").Value = 1545
$ws.Range("J122").Value = 3499.25
$ws.Range("K122").Value = 4635
$ws.Range("L122").Value = 10497.75
$ws.Range("M122").Value = -2185
$ws.Range("N122").Value = -15397.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 858.2
$ws.Range("I64").Value = 796.3333
$ws.Range("K64").Value = 796.3333
$ws.Range("M64").Value = -571.3333
$ws.Range("H67").Value = 858.2
$ws.Range("I67").Value = 796.3333
$ws.Range("K67").Value = 796.3333
$ws.Range("M67").Value = -16.33330000000001
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 1633
$ws.Range("I134").Value = 1416.25
$ws.Range("K134").Value = 4248.75
$ws.Range("M134").Value = -1713.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1500
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H86").Value = 12848
$ws.Range("I86").Value = 5700
$ws.Range("K86").Value = 5700
$ws.Range("M86").Value = -4577
$ws.Range("H89").Value = 12848
$ws.Range("I89").Value = 5700
$ws.Range("K89").Value = 28500
$ws.Range("M89").Value = -22884

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 161.57143
$ws.Range("I17").Value = 220.25
$ws.Range("J17").Value = 83.333336
$ws.Range("K17").Value = 660.75
$ws.Range("L17").Value = 250.000008
$ws.Range("M17").Value = -491.75
$ws.Range("N17").Value = -588.000008
$ws.Range("H109").Value = 2014.25
$ws.Range("I109").Value = 2014.25
$ws.Range("K109").Value = 6042.75
$ws.Range("M109").Value = -5002.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2446
$ws.Range("H107").Value = 1462.9
$ws.Range("I107").Value = 564.2
$ws.Range("J107").Value = 2361.6
$ws.Range("K107").Value = 564.2
$ws.Range("L107").Value = 2361.6
$ws.Range("M107").Value = 1355.8
$ws.Range("N107").Value = -6201.6
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3750.5334
$ws.Range("I40").Value = 3446.2
$ws.Range("K40").Value = 3446.2
$ws.Range("M40").Value = -3310.2
$ws.Range("H46").Value = 2728.4443
$ws.Range("I46").Value = 2096.5
$ws.Range("J46").Value = 3234
$ws.Range("K46").Value = 2096.5
$ws.Range("L46").Value = 3234
$ws.Range("M46").Value = -1908.5
$ws.Range("N46").Value = -3610
$ws.Range("H93").Value = 1609
$ws.Range("I93").Value = 1595.5
$ws.Range("J93").Value = 1649.5
$ws.Range("K93").Value = 1595.5
$ws.Range("L93").Value = 1649.5
$ws.Range("M93").Value = -347.5
$ws.Range("N93").Value = -4145.5
$ws.Range("H100").Value = 5130.1
$ws.Range("I100").Value = 5130.1
$ws.Range("K100").Value = 5130.1
$ws.Range("M100").Value = -4589.1

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1112222.5
$ws.Range("I81").Value = 1250.125
$ws.Range("K81").Value = 2500.25
$ws.Range("M81").Value = -1439.25
$ws.Range("H84").Value = 1112222.5
$ws.Range("I84").Value = 1250.125
$ws.Range("K84").Value = 12501.25
$ws.Range("M84").Value = -7197.25
$ws.Range("H122").Value = 6787.6665
$ws.Range("I122").Value = 6080.5713
$ws.Range("K122").Value = 18241.7139
$ws.Range("M122").Value = -15791.7139
$ws.Range("H126").Value = 1521.5555
$ws.Range("I126").Value = 1128
$ws.Range("K126").Value = 3384
$ws.Range("M126").Value = -914
$ws.Range("H136").Value = 3182.3333
$ws.Range("I136").Value = 3182.3333
$ws.Range("K136").Value = 9546.999899999999
$ws.Range("M136").Value = -6996.999899999999
